$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gastos")

$ws.Range("A1").Value = "Descripción"
$ws.Range("B1").Value = "Monto"
$ws.Range("C1").Value = "Categoría"
$ws.Range("D1").Value = "Nivel de necesidad"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 5
